$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Russia/successor states to soviet union" row (code 16, row 8) is being
# retired as its own entry: its concept is folded into the "Europe, Central
# Asia, Transcaucasia" row (row 7), which becomes "Europe, Russia, Central
# Asia, Transcaucasia". So we update row 7's text first, then delete row 8
# entirely (Excel will shift everything below it up by one row).

$ws.Cells.Item(7, 2).Value = "Europe, Russia, Central Asia, Transcaucasia"
$ws.Cells.Item(7, 3).Value = "Combine Europe, Other European Countries, Russia, Baltic States, Central Asia and Transcaucasia - Baltic states include Estonia, Latvia and Lithuania | Central Asia includes Kazakhstan, Krygyztan, Tajikistan, Turkmenistan, and Uzbekistan | Transcaucasia includes Armenia, Azerbaijan, and Georgia"
$ws.Rows.Item(7).RowHeight = 85.5

$ws.Rows.Item(8).Delete()

# Move the "Combine India, Pakistan, and Bangladesh" note up one row, from the
# India row to the Korea row - and clear it off of the India row (net effect
# after the deletion above: it now lives on C10 instead of C11).
$ws.Cells.Item(10, 3).Value = "Combine India, Pakistan, and Bangladesh"
$ws.Cells.Item(11, 3).ClearContents()

# Reflect the new selection recorded in the sheet view.
$ws.Range("C10:C11").Select()
